$d = $word.ActiveDocument

function Find-ParaIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# Locate every paragraph we need to touch up front (positions shift once we
# start deleting, so we resolve indices first and then work from the bottom
# of the document upward so earlier indices stay valid).
$howIdx    = Find-ParaIndex("^How It Works in Your Project")
$outputIdx = Find-ParaIndex("^Output: A trained LSTM model")
$hrIdx     = $howIdx - 1          # paragraph that only held the <w:pict> hr
$djIdx     = Find-ParaIndex("^Django")
$pgIdx     = Find-ParaIndex("^PostgreSQL")

# ---------------------------------------------------------------------------
# 4) Remove the old _GoBack bookmark sitting in front of "PostgreSQL" (it is
#    relocated to the horizontal-rule paragraph below). Done first since it
#    is the last affected paragraph in the document.
# ---------------------------------------------------------------------------
$pPg = $d.Paragraphs($pgIdx)
$pgFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w:rsidR="0022120C" w:rsidRDefault="009C68E5" w:rsidP="009C68E5">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr>' +
    '<w:r w:rsidRPr="009C68E5"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>PostgreSQL</w:t></w:r>' +
    '<w:r w:rsidRPr="009C68E5"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>: A reliable, scalable, and feature-rich database with support for advanced querying, data integrity, and handling large-scale applications efficiently.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$pPg.Range.InsertXML($pgFrag)

# ---------------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> in front of "Django" (it was
#    only valid while the now-removed content pushed it onto a new page).
# ---------------------------------------------------------------------------
$pDj = $d.Paragraphs($djIdx)
$djFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w:rsidR="009C68E5" w:rsidRPr="009C68E5" w:rsidRDefault="009C68E5" w:rsidP="009C68E5">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
      '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="009C68E5"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Django</w:t></w:r>' +
    '<w:r w:rsidRPr="009C68E5"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>: Ideal for building secure, scalable, and fast backend systems with its powerful ORM, built-in security features, and rapid development capabilities.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$pDj.Range.InsertXML($djFrag)

# ---------------------------------------------------------------------------
# 1) Remove the "How It Works in Your Project" ... "Short Summary" block
#    (from the heading through the "Output: ..." bullet, inclusive), which
#    sits between the horizontal-rule paragraph and the "Django" bullet.
# ---------------------------------------------------------------------------
$pHow = $d.Paragraphs($howIdx)
$pOutput = $d.Paragraphs($outputIdx)
$delRange = $d.Range($pHow.Range.Start, $pOutput.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 2) The paragraph that used to hold only the horizontal-rule <w:pict> is now
#    empty; clear it out and stamp it with the relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$pHr = $d.Paragraphs($hrIdx)
$hrFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w:rsidR="009C68E5" w:rsidRPr="009C68E5" w:rsidRDefault="009C68E5" w:rsidP="009C68E5">' +
    '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
      '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '</w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$pHr.Range.InsertXML($hrFrag)

Write-Host "Edit applied."
